$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Tgfb1", "Sdc2", "ECs", 3, 1, 51.59157666666666, 154.77473, 0.2641250550177587, 0.2641250550177588, 3, 1, 3.057109, 9.171327, 0.02694952608666365, 0.02694952608666365, 157.7210733518567, 1419.48966016671, 0.007118045060342561, 0.007118045060342562),
    @("ECs", "Tgfb1", "Sdc2", "FAPs", 3, 1, 51.59157666666666, 154.77473, 0.2641250550177587, 0.2641250550177588, 3, 1, 89.02756, 267.08268, 0.7848102735793893, 0.7848102735793893, 4593.072187186266, 41337.64968467639, 0.2072880566876585, 0.2072880566876585),
    @("ECs", "Tgfb1", "Sdc2", "M2", 3, 1, 51.59157666666666, 154.77473, 0.2641250550177587, 0.2641250550177588, 2, 0.6666666666666666, 0.184005, 0.5520149999999999, 0.0016220709001794, 0.0016220709001794, 9.493108064549999, 85.43797258094997, 0.0004284295657525894, 0.0004284295657525894),
    @("ECs", "Tgfb1", "Sdc2", "sCs", 3, 1, 51.59157666666666, 154.77473, 0.2641250550177587, 0.2641250550177588, 3, 1, 21.16964733333333, 63.508942, 0.1866181294337677, 0.1866181294337677, 1092.175483403962, 9829.579350635659, 0.04929052370400513, 0.04929052370400513),
    @("FAPs", "Tgfb1", "Sdc2", "ECs", 3, 1, 19.32115333333334, 57.96346000000001, 0.09891538535728452, 0.09891538535728453, 3, 1, 3.057109, 9.171327, 0.02694952608666365, 0.02694952608666365, 59.06687174571335, 531.6018457114201, 0.002665722758058527, 0.002665722758058527),
    @("FAPs", "Tgfb1", "Sdc2", "FAPs", 3, 1, 19.32115333333334, 57.96346000000001, 0.09891538535728452, 0.09891538535728453, 3, 1, 89.02756, 267.08268, 0.7848102735793893, 0.7848102735793893, 1720.115137652534, 15481.0362388728, 0.07762981064346118, 0.07762981064346118),
    @("FAPs", "Tgfb1", "Sdc2", "M2", 3, 1, 19.32115333333334, 57.96346000000001, 0.09891538535728452, 0.09891538535728453, 2, 0.6666666666666666, 0.184005, 0.5520149999999999, 0.0016220709001794, 0.0016220709001794, 3.555188819100001, 31.9966993719, 0.0001604477681680827, 0.0001604477681680827),
    @("FAPs", "Tgfb1", "Sdc2", "sCs", 3, 1, 19.32115333333334, 57.96346000000001, 0.09891538535728452, 0.09891538535728453, 3, 1, 21.16964733333333, 63.508942, 0.1866181294337677, 0.1866181294337677, 409.0220021399246, 3681.198019259321, 0.01845940418759674, 0.01845940418759674),
    @("M2", "Tgfb1", "Sdc2", "ECs", 3, 1, 112.3724673333333, 337.117402, 0.5752951554216499, 0.57529515542165, 3, 1, 3.057109, 9.171327, 0.02694952608666365, 0.02694952608666365, 343.5348812369393, 3091.813931132453, 0.01550393179856697, 0.01550393179856698),
    @("M2", "Tgfb1", "Sdc2", "FAPs", 3, 1, 112.3724673333333, 337.117402, 0.5752951554216499, 0.57529515542165, 3, 1, 89.02756, 267.08268, 0.7848102735793893, 0.7848102735793893, 10004.24657786637, 90038.21920079735, 0.4514975483153623, 0.4514975483153624),
    @("M2", "Tgfb1", "Sdc2", "M2", 3, 1, 112.3724673333333, 337.117402, 0.5752951554216499, 0.57529515542165, 2, 0.6666666666666666, 0.184005, 0.5520149999999999, 0.0016220709001794, 0.0016220709001794, 20.67709585166999, 186.0938626650299, 0.0009331695306236433, 0.0009331695306236434),
    @("M2", "Tgfb1", "Sdc2", "sCs", 3, 1, 112.3724673333333, 337.117402, 0.5752951554216499, 0.57529515542165, 3, 1, 21.16964733333333, 63.508942, 0.1866181294337677, 0.1866181294337677, 2378.885503423187, 21409.96953080868, 0.107360505777097, 0.107360505777097),
    @("sCs", "Tgfb1", "Sdc2", "ECs", 3, 1, 12.044915, 36.134745, 0.06166440420330686, 0.06166440420330688, 3, 1, 3.057109, 9.171327, 0.02694952608666365, 0.02694952608666365, 36.822618050735, 331.403562456615, 0.00166182646969559, 0.00166182646969559),
    @("sCs", "Tgfb1", "Sdc2", "FAPs", 3, 1, 12.044915, 36.134745, 0.06166440420330686, 0.06166440420330688, 3, 1, 89.02756, 267.08268, 0.7848102735793893, 0.7848102735793893, 1072.3293928574, 9650.9645357166, 0.0483948579329073, 0.04839485793290731),
    @("sCs", "Tgfb1", "Sdc2", "M2", 3, 1, 12.044915, 36.134745, 0.06166440420330686, 0.06166440420330688, 2, 0.6666666666666666, 0.184005, 0.5520149999999999, 0.0016220709001794, 0.0016220709001794, 2.216324584575, 19.946921261175, 0.0001000240356350843, 0.0001000240356350843),
    @("sCs", "Tgfb1", "Sdc2", "sCs", 3, 1, 12.044915, 36.134745, 0.06166440420330686, 0.06166440420330688, 3, 1, 21.16964733333333, 63.508942, 0.1866181294337677, 0.1866181294337677, 254.9866027099767, 2294.87942438979, 0.01150769576506889, 0.01150769576506889)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
